# Changed product scroll to uniform grid: widen the THEME column and
# append a new product theme / foreground color entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (THEME) so the longer theme names fit in the uniform grid.
# (ColumnWidth snaps to the nearest pixel boundary, 23.59 is the closest
# input that lands on the intended ~24.4 character width.)
$ws.Columns.Item(3).ColumnWidth = 23.59

# Add the new theme / foreground color pair as row 28.
$ws.Range("C28").Value = "btnDarkYellowItemTheme"
$ws.Range("D28").Value = "Red"
